# Auto-generated PowerShell Word COM-interop script
# Replaces division-problem answer strings in table cells per the diff.
$d = $word.ActiveDocument

$d.Content.Find.Execute("746÷8=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "374÷2=187, 0", 2) | Out-Null
$d.Content.Find.Execute("985÷3=328, 1", $true, $false, $false, $false, $false, $true, 1, $false, "161÷3=53, 2", 2) | Out-Null
$d.Content.Find.Execute("893÷8=111, 5", $true, $false, $false, $false, $false, $true, 1, $false, "213÷6=35, 3", 2) | Out-Null
$d.Content.Find.Execute("445÷9=49, 4", $true, $false, $false, $false, $false, $true, 1, $false, "846÷9=94, 0", 2) | Out-Null
$d.Content.Find.Execute("612÷8=76, 4", $true, $false, $false, $false, $false, $true, 1, $false, "297÷5=59, 2", 2) | Out-Null
$d.Content.Find.Execute("272÷9=30, 2", $true, $false, $false, $false, $false, $true, 1, $false, "337÷5=67, 2", 2) | Out-Null
$d.Content.Find.Execute("588÷2=294, 0", $true, $false, $false, $false, $false, $true, 1, $false, "664÷4=166, 0", 2) | Out-Null
$d.Content.Find.Execute("385÷5=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "756÷9=84, 0", 2) | Out-Null
$d.Content.Find.Execute("276÷6=46, 0", $true, $false, $false, $false, $false, $true, 1, $false, "637÷4=159, 1", 2) | Out-Null
$d.Content.Find.Execute("303÷2=151, 1", $true, $false, $false, $false, $false, $true, 1, $false, "407÷2=203, 1", 2) | Out-Null
$d.Content.Find.Execute("550÷3=183, 1", $true, $false, $false, $false, $false, $true, 1, $false, "722÷3=240, 2", 2) | Out-Null
$d.Content.Find.Execute("353÷3=117, 2", $true, $false, $false, $false, $false, $true, 1, $false, "390÷7=55, 5", 2) | Out-Null
$d.Content.Find.Execute("137÷2=68, 1", $true, $false, $false, $false, $false, $true, 1, $false, "227÷6=37, 5", 2) | Out-Null
$d.Content.Find.Execute("851÷7=121, 4", $true, $false, $false, $false, $false, $true, 1, $false, "198÷6=33, 0", 2) | Out-Null
$d.Content.Find.Execute("649÷8=81, 1", $true, $false, $false, $false, $false, $true, 1, $false, "242÷9=26, 8", 2) | Out-Null
$d.Content.Find.Execute("183÷9=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "266÷9=29, 5", 2) | Out-Null
$d.Content.Find.Execute("916÷6=152, 4", $true, $false, $false, $false, $false, $true, 1, $false, "134÷9=14, 8", 2) | Out-Null
$d.Content.Find.Execute("386÷8=48, 2", $true, $false, $false, $false, $false, $true, 1, $false, "194÷9=21, 5", 2) | Out-Null
$d.Content.Find.Execute("110÷9=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "793÷7=113, 2", 2) | Out-Null
$d.Content.Find.Execute("126÷5=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "730÷2=365, 0", 2) | Out-Null
$d.Content.Find.Execute("688÷2=344, 0", $true, $false, $false, $false, $false, $true, 1, $false, "255÷2=127, 1", 2) | Out-Null
$d.Content.Find.Execute("750÷7=107, 1", $true, $false, $false, $false, $false, $true, 1, $false, "826÷9=91, 7", 2) | Out-Null
$d.Content.Find.Execute("876÷6=146, 0", $true, $false, $false, $false, $false, $true, 1, $false, "897÷9=99, 6", 2) | Out-Null
$d.Content.Find.Execute("508÷9=56, 4", $true, $false, $false, $false, $false, $true, 1, $false, "529÷3=176, 1", 2) | Out-Null
$d.Content.Find.Execute("216÷6=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "173÷4=43, 1", 2) | Out-Null
